$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-18 Saturday" "2025-10-19 Sunday"

Replace-Text "496÷2=248, 0" "741÷6=123, 3"
Replace-Text "253÷4=63, 1" "145÷6=24, 1"
Replace-Text "283÷6=47, 1" "744÷4=186, 0"
Replace-Text "244÷7=34, 6" "281÷4=70, 1"
Replace-Text "993÷5=198, 3" "205÷2=102, 1"

Replace-Text "577÷3=192, 1" "324÷5=64, 4"
Replace-Text "753÷8=94, 1" "973÷5=194, 3"
Replace-Text "337÷9=37, 4" "434÷7=62, 0"
Replace-Text "440÷8=55, 0" "294÷9=32, 6"
Replace-Text "741÷5=148, 1" "394÷3=131, 1"

Replace-Text "646÷6=107, 4" "372÷6=62, 0"
Replace-Text "131÷2=65, 1" "768÷7=109, 5"
Replace-Text "371÷6=61, 5" "477÷5=95, 2"
Replace-Text "581÷4=145, 1" "362÷8=45, 2"
Replace-Text "340÷8=42, 4" "304÷5=60, 4"

Replace-Text "148÷5=29, 3" "866÷4=216, 2"
Replace-Text "534÷6=89, 0" "423÷6=70, 3"
Replace-Text "225÷7=32, 1" "798÷2=399, 0"
Replace-Text "439÷8=54, 7" "138÷6=23, 0"
Replace-Text "737÷7=105, 2" "489÷6=81, 3"

Replace-Text "993÷3=331, 0" "597÷8=74, 5"
Replace-Text "621÷8=77, 5" "119÷3=39, 2"
Replace-Text "703÷3=234, 1" "851÷5=170, 1"
Replace-Text "566÷4=141, 2" "224÷2=112, 0"
Replace-Text "234÷2=117, 0" "887÷4=221, 3"

Write-Output "Done"
